# Execution of DemoWebShop testcases inprogress
# Fill in the next set of generated test data rows for the DemoWebShop
# registration / address / multiple-products test cases.

$wb = $excel.ActiveWorkbook

# --- Demo_Registration: new registration test data (row 2) ---
$wsReg = $wb.Worksheets.Item("Demo_Registration")
$wsReg.Range("C2").Value = "Allan"
$wsReg.Range("D2").Value = "Luettgen"
$wsReg.Range("E2").Value = "QhdEmipn@test.org"

# --- Demo_CreateAddress: new address test data (row 2) ---
$wsAddr = $wb.Worksheets.Item("Demo_CreateAddress")
$wsAddr.Range("E2").Value = "Nickole"
$wsAddr.Range("F2").Value = "Schultz"
$wsAddr.Range("G2").Value = "fvrfvahy@test.org"
$wsAddr.Range("H2").Value = "Schneider Group"
$wsAddr.Range("I2").Value = "Mozambique"
$wsAddr.Range("K2").Value = "New Marcustown"
$wsAddr.Range("L2").Value = "'189"
$wsAddr.Range("M2").Value = "621 Haag Green"
$wsAddr.Range("N2").Value = "DM"
$wsAddr.Range("O2").Value = "910-797-3373"

# --- Demo_MultipleProducts: new order numbers (M2:M6) ---
$wsOrders = $wb.Worksheets.Item("Demo_MultipleProducts")
$wsOrders.Range("M2").Value = "Order number: 1576750"
$wsOrders.Range("M3").Value = "Order number: 1576753"
$wsOrders.Range("M4").Value = "Order number: 1576757"
$wsOrders.Range("M5").Value = "Order number: 1576766"
$wsOrders.Range("M6").Value = "Order number: 1576770"
